$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09102500000000001
$ws.Range("H2").Value = 0.273075
$ws.Range("I2").Value = 0.8515285885346505
$ws.Range("J2").Value = 0.8515285885346504
$ws.Range("M2").Value = 0.2557903333333333
$ws.Range("N2").Value = 0.767371
$ws.Range("O2").Value = 0.1295099616231651
$ws.Range("P2").Value = 0.1295099616231651
$ws.Range("Q2").Value = 0.02328331509166667
$ws.Range("R2").Value = 0.209549835825
$ws.Range("S2").Value = 0.1102814348221505
$ws.Range("T2").Value = 0.1102814348221505
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09102500000000001
$ws.Range("H3").Value = 0.273075
$ws.Range("I3").Value = 0.8515285885346505
$ws.Range("J3").Value = 0.8515285885346504
$ws.Range("O3").Value = 0.7320547918387076
$ws.Range("P3").Value = 0.7320547918387076
$ws.Range("Q3").Value = 0.1316088906916667
$ws.Range("R3").Value = 1.184480016225
$ws.Range("S3").Value = 0.6233655836244421
$ws.Range("T3").Value = 0.6233655836244419
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09102500000000001
$ws.Range("H4").Value = 0.273075
$ws.Range("I4").Value = 0.8515285885346505
$ws.Range("J4").Value = 0.8515285885346504
$ws.Range("O4").Value = 0.1384352465381273
$ws.Range("P4").Value = 0.1384352465381273
$ws.Range("Q4").Value = 0.02488790379166667
$ws.Range("R4").Value = 0.223991134125
$ws.Range("S4").Value = 0.1178815700880579
$ws.Range("T4").Value = 0.1178815700880579
$ws.Range("I5").Value = 0.1484714114653495
$ws.Range("J5").Value = 0.1484714114653495
$ws.Range("M5").Value = 0.2557903333333333
$ws.Range("N5").Value = 0.767371
$ws.Range("O5").Value = 0.1295099616231651
$ws.Range("P5").Value = 0.1295099616231651
$ws.Range("Q5").Value = 0.004059648380333333
$ws.Range("R5").Value = 0.036536835423
$ws.Range("S5").Value = 0.01922852680101457
$ws.Range("T5").Value = 0.01922852680101457
$ws.Range("I6").Value = 0.1484714114653495
$ws.Range("J6").Value = 0.1484714114653495
$ws.Range("O6").Value = 0.7320547918387076
$ws.Range("P6").Value = 0.7320547918387076
$ws.Range("S6").Value = 0.1086892082142655
$ws.Range("T6").Value = 0.1086892082142655
$ws.Range("I7").Value = 0.1484714114653495
$ws.Range("J7").Value = 0.1484714114653495
$ws.Range("O7").Value = 0.1384352465381273
$ws.Range("P7").Value = 0.1384352465381273
$ws.Range("S7").Value = 0.02055367645006939
$ws.Range("T7").Value = 0.0205536764500694
